$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.834.53"
$ws.Range("E2").Value = "  -2.01%  "

$ws.Range("D3").Value = "2.277.21"
$ws.Range("E3").Value = "  -3.36%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.62%  "

$ws.Range("E7").Value = "  -1.41%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.603"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.98%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.83"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.65%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0905"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.39%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.23"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.52%  "

$ws.Range("E13").Value = "  -0.37%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.955"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.77%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.22"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.38%  "

$ws.Range("D16").Value = "2.621.76"
$ws.Range("E16").Value = "  -3.09%  "

$ws.Range("D17").Value = "2.281.31"
$ws.Range("E17").Value = "  -2.68%  "

$ws.Range("D18").Value = "41.684.24"
$ws.Range("E18").Value = "  -2.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.40%  "

$ws.Range("E20").Value = "  -1.39%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "280.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +9.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.80%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.95%  "

$ws.Range("E25").Value = "  -3.60%  "

$ws.Range("E26").Value = "  +0.79%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.35%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.13%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.25%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "163.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.84%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.97%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0870"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.52%  "

$ws.Range("E33").Value = "  -0.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.23%  "

$ws.Range("E35").Value = "  +0.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.114"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -9.41%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.55"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.54%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.88"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.32%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0346"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.89%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.63"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "101.88"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +18.24%  "

$ws.Range("E42").Value = "  -2.52%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "69.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.11%  "

$ws.Range("E44").Value = "  -0.12%  "

$ws.Range("E45").Value = "  -6.40%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "115.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.38%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.97%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "75.49"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.00%  "

$ws.Range("E50").Value = "  -4.53%  "

$ws.Range("E51").Value = "  -2.72%  "
